$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: new person "Debbie Tarenskeen" ---
$ws.Range("A7").Value = "Debbie Tarenskeen"
$ws.Range("B7").Value = "Debbie"
$ws.Range("D7").Value = "Tarenskeen"

# --- Row 18: new account for Debbie at HAN ---
$ws.Range("A18").Value = "Acc_Debbie"
$ws.Range("B18").Value = "Debbie"
$ws.Range("C18").Value = "welkom"
$ws.Range("D18").Value = "Debbie Tarenskeen"
$ws.Range("E18").Value = "HAN"
$ws.Range("F18").Value = "Tutor"
$ws.Range("G18").Value = "AccountMgr"
$ws.Range("H18").Value = "Student"

# --- Row 11: new organization "HAN" ---
$ws.Range("A11").Value = "HAN"
$ws.Range("B11").Value = "HAN"
$ws.Range("C11").Value = "Hogeschool Arnhem Nijmegen"

# --- Selection moved to C12 ---
$ws.Range("C12").Select()
